$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.086.12"
$ws.Range("E2").Value = "  -4.40%  "
$ws.Range("D3").Value = "2.619.08"
$ws.Range("E3").Value = "  -3.06%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.16"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.75"
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.565"
$ws.Range("E8").Value = "  -2.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.67"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("E10").Value = "  -3.28%  "
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("D13").Value = "3.079.31"
$ws.Range("E13").Value = "  -3.17%  "
$ws.Range("D14").Value = "58.086.94"
$ws.Range("E14").Value = "  -4.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.63"
$ws.Range("E15").Value = "  -3.23%  "
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("D17").Value = "2.614.35"
$ws.Range("E17").Value = "  -8.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.38"
$ws.Range("E18").Value = "  -3.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "334.23"
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.34"
$ws.Range("E20").Value = "  -2.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.25"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.15"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.165"
$ws.Range("E25").Value = "  -2.76%  "
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.03"
$ws.Range("E27").Value = "  -3.53%  "
$ws.Range("D28").Value = "0.0₃0780"
$ws.Range("E28").Value = "  -4.82%  "
$ws.Range("E29").Value = "  -2.94%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.58"
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.86"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.66"
$ws.Range("E33").Value = "  -2.22%  "
$ws.Range("E34").Value = "  -4.25%  "
$ws.Range("E35").Value = "  -5.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.896"
$ws.Range("E36").Value = "  -4.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.61"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.844"
$ws.Range("E38").Value = "  -3.48%  "
$ws.Range("E39").Value = "  -6.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.59"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("E42").Value = "  -1.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0964"
$ws.Range("E43").Value = "  -2.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "267.97"
$ws.Range("E44").Value = "  -5.36%  "
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.08"
$ws.Range("E46").Value = "  -4.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0531"
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("D48").Value = "2.027.86"
$ws.Range("E48").Value = "  -5.30%  "
$ws.Range("E49").Value = "  -2.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.63"
$ws.Range("E50").Value = "  -5.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.18"
$ws.Range("E51").Value = "  -5.12%  "
